# Apply updated stimuli order to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 96, "dog/dog100.png", "prüfen", "dog"),
    @(3, 126, "dog/dog069.png", "danken", "dog"),
    @(4, 6, "dog/dog088.png", "wehen", "dog"),
    @(5, 76, "car/car108.png", "deuten", "car"),
    @(6, 71, "car/car107.png", "quellen", "car"),
    @(7, 2, "car/car095.png", "bauen", "car"),
    @(8, 49, "dog/dog104.png", "kennen", "dog"),
    @(9, 53, "car/car078.png", "bergen", "car"),
    @(10, 3, "car/car072.png", "ändern", "car"),
    @(11, 65, "dog/dog110.png", "sparen", "dog"),
    @(12, 61, "car/car084.png", "trotzen", "car"),
    @(13, 44, "car/car071.png", "narren", "car"),
    @(14, 29, "dog/dog082.png", "proben", "dog"),
    @(15, 79, "car/car067.png", "atmen", "car"),
    @(16, 77, "dog/dog066.png", "stoßen", "dog"),
    @(17, 43, "car/car073.png", "spüren", "car"),
    @(18, 30, "car/car091.png", "dienen", "car"),
    @(19, 11, "dog/dog070.png", "leeren", "dog"),
    @(20, 60, "dog/dog084.png", "passen", "dog"),
    @(21, 100, "car/car093.png", "streifen", "car"),
    @(22, 14, "dog/dog087.png", "herrschen", "dog"),
    @(23, 127, "dog/dog073.png", "wecken", "dog"),
    @(24, 92, "car/car111.png", "münzen", "car"),
    @(25, 84, "dog/dog101.png", "frischen", "dog"),
    @(26, 122, "dog/dog092.png", "nullen", "dog"),
    @(27, 120, "car/car122.png", "süßen", "car"),
    @(28, 111, "dog/dog095.png", "grenzen", "dog"),
    @(29, 41, "car/car075.png", "spenden", "car"),
    @(30, 0, "dog/dog081.png", "meinen", "dog"),
    @(31, 12, "car/car094.png", "holen", "car"),
    @(32, 62, "dog/dog071.png", "kranken", "dog"),
    @(33, 121, "car/car082.png", "lassen", "car")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
